$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Score values for columns E (Specific) through J (Relevant), rows 2-18
$scores = @{
    2  = @(2,2,1,2,2,2)
    3  = @(2,2,1,2,2,2)
    4  = @(2,1,1,1,2,2)
    5  = @(2,2,1,1,2,2)
    6  = @(2,2,1,1,2,2)
    7  = @(2,2,1,1,2,2)
    8  = @(2,2,1,1,2,2)
    9  = @(2,1,2,1,1,2)
    10 = @(2,2,1,1,2,2)
    11 = @(2,2,1,1,2,2)
    12 = @(2,2,1,1,2,2)
    13 = @(2,1,2,1,2,2)
    14 = @(2,2,1,1,2,2)
    15 = @(2,1,1,1,1,2)
    16 = @(2,2,1,1,2,2)
    17 = @(2,2,1,1,1,2)
    18 = @(2,2,1,1,2,2)
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # E=5 .. J=10
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# View settings: zoom, frozen header row, final selection (matches the
# sheetView produced by the author scrolling to row 14 with the top row
# frozen, then clicking E16).
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E16").Select() | Out-Null
